$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update 2019/2020 (columns P/Q) figures before removing 2021/2022 columns
$ws.Range("P5").Value2 = 27
$ws.Range("Q5").Value2 = 25.3

$ws.Range("P6").Value2 = 19.6
$ws.Range("Q6").Value2 = 17.8

$ws.Range("P8").Value2 = 2.2
$ws.Range("Q8").Value2 = 2

$ws.Range("P9").Value2 = 5.2
$ws.Range("Q9").Value2 = 5.5

# Remove the 2021 (R) and 2022 (S) columns entirely
$ws.Range("R1:S12").EntireColumn.Delete() | Out-Null

# Update the active selection to reflect the edited cell
$ws.Range("N13").Select() | Out-Null
